# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 634
$wsExpo.Range("F3").Value = 2211
$wsExpo.Range("F4").Value = 91
$wsExpo.Range("F5").Value = 13225
$wsExpo.Range("F11").Value = 991
$wsExpo.Range("F12").Value = 13787
$wsExpo.Range("F13").Value = 14414
$wsExpo.Range("F21").Value = 39
$wsExpo.Range("F24").Value = 57
$wsExpo.Range("F25").Value = 5465
$wsExpo.Range("F27").Value = 152
$wsExpo.Range("F28").Value = 336
$wsExpo.Range("F29").Value = 25
$wsExpo.Range("F30").Value = 74

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 634
$wsAll.Range("F3").Value = 2211
$wsAll.Range("F4").Value = 91
$wsAll.Range("F5").Value = 13225
$wsAll.Range("F12").Value = 991
$wsAll.Range("F13").Value = 13787
$wsAll.Range("F14").Value = 14414
$wsAll.Range("F22").Value = 39
$wsAll.Range("F25").Value = 57
$wsAll.Range("F26").Value = 5465
$wsAll.Range("F28").Value = 152
$wsAll.Range("F29").Value = 336
$wsAll.Range("F30").Value = 25
$wsAll.Range("F31").Value = 74
